# Apply the data correction to Sheet1: the longitude in G2 should be negative
# (west longitude), and update the active selection to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the sign of the longitude value in G2
$ws.Range("G2").Value = -78.010910030000005

# Update the active cell / selection to G2
$ws.Activate()
$ws.Range("G2").Select()
